$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.037.04"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.471.72"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'577.07"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'146.59"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "2.470.45"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "'28.95"
$ws.Range("E14").Value = "  +6.83%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "2.919.25"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "62.977.58"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "2.472.69"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  +3.59%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "'329.59"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'2.23"
$ws.Range("E22").Value = "  +8.82%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'66.29"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "'666.07"
$ws.Range("E26").Value = "  +6.97%  "
$ws.Range("E27").Value = "  +14.76%  "
$ws.Range("D28").Value = "0.0₃0985"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +767.51%  "
$ws.Range("D31").Value = "'1.45"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("D32").Value = "'8.07"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").Value = "'151.71"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'2.72"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D45").Value = "0.0₆0305"
$ws.Range("E45").Value = "  +6.49%  "
$ws.Range("D46").Value = "'151.28"
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("E47").Value = "  +26.20%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'20.70"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("E51").Value = "  -0.92%  "
